# Commit: "Fruta / hortaliza, semanal" — weekly data refresh.
# A new daily price record for Haba (Vega Central Mapocho de Santiago) is
# inserted as row 138, shifting the existing rows 138-179 down to 139-180.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 138 (pushes old rows 138..179 to 139..180)
$ws.Rows.Item(138).Insert()

# Populate the newly inserted row 138 with the new record
$ws.Cells.Item(138, 1).Value = 9
$ws.Cells.Item(138, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(138, 3).Value = "Metropolitana"
$ws.Cells.Item(138, 4).Value = 44524
$ws.Cells.Item(138, 5).Value = 13
$ws.Cells.Item(138, 6).Value = 100112026
$ws.Cells.Item(138, 7).Value = "Haba"
$ws.Cells.Item(138, 8).Value = "Sin especificar"
$ws.Cells.Item(138, 9).Value = "Primera"
$ws.Cells.Item(138, 10).Value = 52
$ws.Cells.Item(138, 11).Value = 7000
$ws.Cells.Item(138, 12).Value = 8000
$ws.Cells.Item(138, 13).Value = 7500
$ws.Cells.Item(138, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(138, 15).Value = "Carahue"
$ws.Cells.Item(138, 16).Value = 300
$ws.Cells.Item(138, 17).Value = 25
$ws.Cells.Item(138, 18).Value = "Hortaliza"

# Make sure the new date cell uses the same date/time number format as the
# rest of column D.
$ws.Cells.Item(138, 4).NumberFormat = $ws.Cells.Item(139, 4).NumberFormat
